$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(31, 8).Value = 88.75
$ws.Cells.Item(31, 9).Value = 88.75
$ws.Cells.Item(31, 11).Value = 266.25
$ws.Cells.Item(31, 13).Value = -36.25

$ws.Cells.Item(70, 8).Value = 189243
$ws.Cells.Item(70, 10).Value = 189243
$ws.Cells.Item(70, 12).Value = 567729
$ws.Cells.Item(70, 14).Value = -568269

$ws.Cells.Item(73, 8).Value = 189243
$ws.Cells.Item(73, 10).Value = 189243
$ws.Cells.Item(73, 12).Value = 567729
$ws.Cells.Item(73, 14).Value = -569601

$ws.Cells.Item(94, 8).Value = 19999.666
$ws.Cells.Item(94, 9).Value = 19999.666
$ws.Cells.Item(94, 11).Value = 19999.666
$ws.Cells.Item(94, 13).Value = -19548.666

$ws.Cells.Item(111, 8).Value = 4256.75
$ws.Cells.Item(111, 9).Value = 4256.75
$ws.Cells.Item(111, 11).Value = 12770.25
$ws.Cells.Item(111, 13).Value = -9703.25

$ws.Cells.Item(137, 8).Value = 2706.35
$ws.Cells.Item(137, 9).Value = 2059.2144
$ws.Cells.Item(137, 10).Value = 4216.3335
$ws.Cells.Item(137, 11).Value = 6177.6432
$ws.Cells.Item(137, 12).Value = 12649.0005
$ws.Cells.Item(137, 13).Value = -3627.6432
$ws.Cells.Item(137, 14).Value = -17749.0005

$ws.Cells.Item(138, 8).Value = 4068.9375
$ws.Cells.Item(138, 10).Value = 4257.25
$ws.Cells.Item(138, 12).Value = 12771.75
$ws.Cells.Item(138, 14).Value = -23051.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1825.1666
$ws.Cells.Item(45, 9).Value = 1825.1666
$ws.Cells.Item(45, 11).Value = 1825.1666
$ws.Cells.Item(45, 13).Value = -1448.1666

$ws.Cells.Item(61, 8).Value = 1318
$ws.Cells.Item(61, 9).Value = 1338.8182
$ws.Cells.Item(61, 11).Value = 1338.8182
$ws.Cells.Item(61, 13).Value = -1126.8182

$ws.Cells.Item(74, 8).Value = 2308.077
$ws.Cells.Item(74, 9).Value = 1013
$ws.Cells.Item(74, 10).Value = 5222
$ws.Cells.Item(74, 11).Value = 1013
$ws.Cells.Item(74, 12).Value = 5222
$ws.Cells.Item(74, 13).Value = -139
$ws.Cells.Item(74, 14).Value = -6970

$ws.Cells.Item(77, 8).Value = 2308.077
$ws.Cells.Item(77, 9).Value = 1013
$ws.Cells.Item(77, 10).Value = 5222
$ws.Cells.Item(77, 11).Value = 5065
$ws.Cells.Item(77, 12).Value = 26110
$ws.Cells.Item(77, 13).Value = -697
$ws.Cells.Item(77, 14).Value = -34846

$ws.Cells.Item(110, 8).Value = 3641
$ws.Cells.Item(110, 10).Value = 379.5
$ws.Cells.Item(110, 12).Value = 379.5
$ws.Cells.Item(110, 14).Value = -4469.5

$ws.Cells.Item(122, 8).Value = 670419
$ws.Cells.Item(122, 9).Value = 1114087.4
$ws.Cells.Item(122, 10).Value = 4916.5
$ws.Cells.Item(122, 11).Value = 3342262.2
$ws.Cells.Item(122, 12).Value = 14749.5
$ws.Cells.Item(122, 13).Value = -3339812.2
$ws.Cells.Item(122, 14).Value = -19649.5

$ws.Cells.Item(132, 8).Value = 1431.8889
$ws.Cells.Item(132, 9).Value = 1362.762
$ws.Cells.Item(132, 10).Value = 2399.6667
$ws.Cells.Item(132, 11).Value = 4088.286
$ws.Cells.Item(132, 12).Value = 7199.000100000001
$ws.Cells.Item(132, 13).Value = -1558.286
$ws.Cells.Item(132, 14).Value = -12259.0001

$ws.Cells.Item(136, 8).Value = 1318
$ws.Cells.Item(136, 9).Value = 1338.8182
$ws.Cells.Item(136, 11).Value = 4016.4546
$ws.Cells.Item(136, 13).Value = -1466.4546

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 9484.857
$ws.Cells.Item(86, 9).Value = 7848.75
$ws.Cells.Item(86, 11).Value = 7848.75
$ws.Cells.Item(86, 13).Value = -6725.75

$ws.Cells.Item(89, 8).Value = 9484.857
$ws.Cells.Item(89, 9).Value = 7848.75
$ws.Cells.Item(89, 11).Value = 39243.75
$ws.Cells.Item(89, 13).Value = -33627.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 125.666664
$ws.Cells.Item(2, 9).Value = 110.75
$ws.Cells.Item(2, 10).Value = 155.5
$ws.Cells.Item(2, 11).Value = 664.5
$ws.Cells.Item(2, 12).Value = 933
$ws.Cells.Item(2, 13).Value = -551.5
$ws.Cells.Item(2, 14).Value = -1159

$ws.Cells.Item(118, 8).Value = 1499.6666
$ws.Cells.Item(118, 9).Value = 1450
$ws.Cells.Item(118, 10).Value = 1599
$ws.Cells.Item(118, 11).Value = 4350
$ws.Cells.Item(118, 12).Value = 4797
$ws.Cells.Item(118, 13).Value = -3107
$ws.Cells.Item(118, 14).Value = -7283

$ws.Cells.Item(128, 8).Value = 0
$ws.Cells.Item(128, 9).Value = 0
$ws.Cells.Item(128, 11).Value = 0
$ws.Cells.Item(128, 13).ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(24, 8).Value = 184346.03
$ws.Cells.Item(24, 9).Value = 2001335.4
$ws.Cells.Item(24, 10).Value = 19165.182
$ws.Cells.Item(24, 11).Value = 2001335.4
$ws.Cells.Item(24, 12).Value = 19165.182
$ws.Cells.Item(24, 13).Value = -2001162.4
$ws.Cells.Item(24, 14).Value = -19511.182

$ws.Cells.Item(29, 8).Value = 2681166.8
$ws.Cells.Item(29, 10).Value = 15050
$ws.Cells.Item(29, 12).Value = 15050
$ws.Cells.Item(29, 14).Value = -15630

$ws.Cells.Item(70, 8).Value = 5333.3335
$ws.Cells.Item(70, 10).Value = 5000
$ws.Cells.Item(70, 12).Value = 5000
$ws.Cells.Item(70, 14).Value = -5540

$ws.Cells.Item(73, 8).Value = 5333.3335
$ws.Cells.Item(73, 10).Value = 5000
$ws.Cells.Item(73, 12).Value = 5000
$ws.Cells.Item(73, 14).Value = -6872

$ws.Cells.Item(97, 8).Value = 965.9375
$ws.Cells.Item(97, 9).Value = 802.1818
$ws.Cells.Item(97, 10).Value = 1326.2
$ws.Cells.Item(97, 11).Value = 802.1818
$ws.Cells.Item(97, 12).Value = 1326.2
$ws.Cells.Item(97, 13).Value = -306.1818
$ws.Cells.Item(97, 14).Value = -2318.2

$ws.Cells.Item(102, 8).Value = 1729.5714
$ws.Cells.Item(102, 9).Value = 432.46155
$ws.Cells.Item(102, 10).Value = 3837.375
$ws.Cells.Item(102, 11).Value = 432.46155
$ws.Cells.Item(102, 12).Value = 3837.375
$ws.Cells.Item(102, 13).Value = 1189.53845
$ws.Cells.Item(102, 14).Value = -7081.375

$ws.Cells.Item(132, 8).Value = 6444
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 10).Value = 6444
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 12).Value = 19332
$ws.Cells.Item(132, 13).ClearContents()
$ws.Cells.Item(132, 14).Value = -24392

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(58, 8).Value = 70094
$ws.Cells.Item(58, 10).Value = 70094
$ws.Cells.Item(58, 12).Value = 70094
$ws.Cells.Item(58, 14).Value = -70710

$ws.Cells.Item(60, 8).Value = 49999.5
$ws.Cells.Item(60, 10).Value = 49999.5
$ws.Cells.Item(60, 12).Value = 49999.5
$ws.Cells.Item(60, 14).Value = -51643.5

$ws.Cells.Item(61, 8).Value = 3000
$ws.Cells.Item(61, 9).Value = 3000
$ws.Cells.Item(61, 11).Value = 3000
$ws.Cells.Item(61, 13).Value = -2708

$ws.Cells.Item(113, 8).Value = 680
$ws.Cells.Item(113, 9).Value = 806.8889
$ws.Cells.Item(113, 10).Value = 299.33334
$ws.Cells.Item(113, 11).Value = 2420.6667
$ws.Cells.Item(113, 12).Value = 898.0000200000001
$ws.Cells.Item(113, 13).Value = -250.6667000000002
$ws.Cells.Item(113, 14).Value = -5238.00002
